# Applies the cryptos list update described in the commit:
# "Updated cryptos list on Sun Nov 17 08:44:10 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = "'" + '90.643.03'
$ws.Range('E2').Value = '  -0.71%  '

# Row 3
$ws.Range('D3').Value = "'" + '3.152.03'
$ws.Range('E3').Value = '  +1.19%  '

# Row 4
$ws.Range('D4').Value = "'" + '1.00'
$ws.Range('E4').Value = '  +0.00%  '

# Row 5
$ws.Range('D5').Value = "'" + '240.63'
$ws.Range('E5').Value = '  +10.14%  '

# Row 6
$ws.Range('D6').Value = "'" + '641.75'
$ws.Range('E6').Value = '  +3.18%  '

# Row 7
$ws.Range('D7').Value = "'" + '1.08'
$ws.Range('E7').Value = '  +9.09%  '

# Row 8
$ws.Range('D8').Value = "'" + '0.364'
$ws.Range('E8').Value = '  -3.83%  '

# Row 10
$ws.Range('D10').Value = "'" + '3.133.92'
$ws.Range('E10').Value = '  +0.76%  '

# Row 11
$ws.Range('D11').Value = "'" + '0.722'
$ws.Range('E11').Value = '  +1.16%  '

# Row 12
$ws.Range('E12').Value = '  +2.77%  '

# Row 13
$ws.Range('D13').Value = "'" + '36.50'
$ws.Range('E13').Value = '  +5.61%  '

# Row 14
$ws.Range('D14').Value = "'" + '0.0000250'
$ws.Range('E14').Value = '  -1.91%  '

# Row 15
$ws.Range('E15').Value = '  +4.28%  '

# Row 16
$ws.Range('D16').Value = "'" + '90.340.91'
$ws.Range('E16').Value = '  -0.81%  '

# Row 17
$ws.Range('D17').Value = "'" + '3.729.92'
$ws.Range('E17').Value = '  +1.32%  '

# Row 18
$ws.Range('B18').Value = 'SuiNetwork'
$ws.Range('C18').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D18').Value = "'" + '3.79'
$ws.Range('E18').Value = '  +0.61%  '

# Row 19
$ws.Range('B19').Value = 'WrappedEther'
$ws.Range('C19').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D19').Value = "'" + '3.122.74'
$ws.Range('E19').Value = '  -0.11%  '

# Row 20
$ws.Range('D20').Value = "'" + '0.0000218'
$ws.Range('E20').Value = '  +0.15%  '

# Row 21
$ws.Range('D21').Value = "'" + '14.48'
$ws.Range('E21').Value = '  +2.94%  '

# Row 22
$ws.Range('D22').Value = "'" + '449.83'
$ws.Range('E22').Value = '  +3.36%  '

# Row 23
$ws.Range('E23').Value = '  +9.13%  '

# Row 24
$ws.Range('E24').Value = '  +2.80%  '

# Row 25
$ws.Range('E25').Value = '  -1.89%  '

# Row 26
$ws.Range('D26').Value = "'" + '90.83'
$ws.Range('E26').Value = '  +5.07%  '

# Row 27
$ws.Range('D27').Value = "'" + '12.48'
$ws.Range('E27').Value = '  +2.26%  '

# Row 28
$ws.Range('D28').Value = "'" + '3.298.45'
$ws.Range('E28').Value = '  +0.75%  '

# Row 29
$ws.Range('E29').Value = '  +0.06%  '

# Row 30
$ws.Range('D30').Value = "'" + '9.76'
$ws.Range('E30').Value = '  +7.46%  '

# Row 31
$ws.Range('D31').Value = "'" + '0.161'
$ws.Range('E31').Value = '  -3.87%  '

# Row 32
$ws.Range('D32').Value = "'" + '27.19'
$ws.Range('E32').Value = '  +15.14%  '

# Row 33
$ws.Range('E33').Value = '  +30.26%  '

# Row 34
$ws.Range('D34').Value = "'" + '3.89'
$ws.Range('E34').Value = '  +3.76%  '

# Row 35
$ws.Range('D35').Value = "'" + '519.66'
$ws.Range('E35').Value = '  -0.92%  '

# Row 36
$ws.Range('D36').Value = "'" + '0.152'
$ws.Range('E36').Value = '  +3.21%  '

# Row 37
$ws.Range('D37').Value = "'" + '7.18'
$ws.Range('E37').Value = '  +0.50%  '

# Row 38
$ws.Range('D38').Value = "'" + '1.94'
$ws.Range('E38').Value = '  +4.62%  '

# Row 39
$ws.Range('E39').Value = '  +0.96%  '

# Row 41
$ws.Range('D41').Value = "'" + '0.423'
$ws.Range('E41').Value = '  +5.01%  '

# Row 42
$ws.Range('B42').Value = 'WhiteBITCoin'
$ws.Range('C42').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D42').Value = "'" + '22.18'
$ws.Range('E42').Value = '  -0.38%  '

# Row 43
$ws.Range('B43').Value = 'Hedera'
$ws.Range('C43').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D43').Value = "'" + '0.0863'
$ws.Range('E43').Value = '  -3.25%  '

# Row 45
$ws.Range('D45').Value = "'" + '3.40'
$ws.Range('E45').Value = '  +46.96%  '

# Row 46
$ws.Range('E46').Value = '  +1.40%  '

# Row 47
$ws.Range('D47').Value = "'" + '0.706'
$ws.Range('E47').Value = '  +13.05%  '

# Row 48
$ws.Range('D48').Value = "'" + '151.35'
$ws.Range('E48').Value = '  +1.84%  '

# Row 49
$ws.Range('D49').Value = "'" + '46.12'
$ws.Range('E49').Value = '  +5.07%  '

# Row 50
$ws.Range('D50').Value = "'" + '4.59'
$ws.Range('E50').Value = '  +8.65%  '

# Row 51
$ws.Range('D51').Value = "'" + '1.36'
$ws.Range('E51').Value = '  +4.06%  '
